# Update of ms figures: flip the "model_forcing_estimate" (column I) flag
# from TRUE to FALSE for four parameter rows, clearing the highlight fill
# that was used to flag them, and move the active selection to I21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(18, 22, 23, 24)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 9)   # column I
    $cell.Value = $false
    $cell.Interior.Pattern = -4142  # xlNone -> clear the highlight fill
}

# Move the selection the way the author left it after editing row 21/22.
$ws.Range("I21").Select()
